$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "#rey-van-bacchanten,-en-menaden"
$ws.Range("C2").Value = "Rey van Bacchanten, en Menaden"
$ws.Range("D2").Value = ""

$ws.Range("B3").Value = "#dichtkunde"
$ws.Range("C3").Value = "Dichtkunde"
$ws.Range("D3").Value = ""

$ws.Range("B4").Value = "#danskunst,-dichtkunde-en-maatzang-te-zaamen"
$ws.Range("C4").Value = "Danskunst, Dichtkunde en Maatzang te zaamen"
$ws.Range("D4").Value = ""

$ws.Range("B5").Value = "#rey-van-alle"
$ws.Range("C5").Value = "Rey van alle"
$ws.Range("D5").Value = ""

$ws.Range("B6").Value = "#een-der-landlieden"
$ws.Range("C6").Value = "Een der Landlieden"
$ws.Range("D6").Value = ""

$ws.Range("B7").Value = "#ceres"
$ws.Range("C7").Value = "Ceres"
$ws.Range("D7").Value = ""

$ws.Range("B8").Value = "#megere"
$ws.Range("C8").Value = "Megere"
$ws.Range("D8").Value = ""

$ws.Range("B9").Value = "#een-van-de-landmannen"
$ws.Range("C9").Value = "Een van de Landmannen"
$ws.Range("D9").Value = ""

$ws.Range("B10").Value = "#deze-der-menaden-nevens-de-2"
$ws.Range("C10").Value = "Deze der Menaden nevens de 2"
$ws.Range("D10").Value = ""

$ws.Range("B11").Value = "#sylenus"
$ws.Range("C11").Value = "Sylenus"
$ws.Range("D11").Value = ""

$ws.Range("B12").Value = "#cupido"
$ws.Range("C12").Value = "Cupido"
$ws.Range("D12").Value = ""

$ws.Range("B13").Value = "#een-uit-het-gevolg-van-bacchus"
$ws.Range("C13").Value = "Een uit het gevolg van Bacchus"
$ws.Range("D13").Value = ""

$ws.Range("B14").Value = "#een-uit-het-gevolg-van-ceres"
$ws.Range("C14").Value = "Een uit het Gevolg van Ceres"
$ws.Range("D14").Value = ""

$ws.Range("B15").Value = "#een-uyt-het-gevolg-van-bacchus,-houdende-een-drinkvat-in-de-hand"
$ws.Range("C15").Value = "Een uyt het gevolg van Bacchus, houdende een drinkvat in de hand"
$ws.Range("D15").Value = ""

$ws.Range("B16").Value = "#eerste-der-landlieden"
$ws.Range("C16").Value = "Eerste der Landlieden"
$ws.Range("D16").Value = ""

$ws.Range("B17").Value = "#bacchus"
$ws.Range("C17").Value = "Bacchus"
$ws.Range("D17").Value = ""

$ws.Range("B18").Value = "#eerste-der-menaden"
$ws.Range("C18").Value = "Eerste der Menaden"
$ws.Range("D18").Value = ""

$ws.Range("B19").Value = "#eenige-uit-de-zanggodinnen,-dichtkunde,-maatzang,-dans--en-speelkunst"
$ws.Range("C19").Value = "Eenige uit de Zanggodinnen, Dichtkunde, Maatzang, Dans- en Speelkunst"
$ws.Range("D19").Value = ""

$ws.Range("B20").Value = "#venus"
$ws.Range("C20").Value = "Venus"
$ws.Range("D20").Value = ""

$ws.Range("B21").Value = "#iris"
$ws.Range("C21").Value = "Iris"
$ws.Range("D21").Value = ""

$ws.Range("B22").Value = "#juno"
$ws.Range("C22").Value = "Juno"
$ws.Range("D22").Value = ""

$ws.Range("B23").Value = "#jupiter"
$ws.Range("C23").Value = "Jupiter"
$ws.Range("D23").Value = ""

$ws.Range("B24").Value = "#tweede-der-landlieden"
$ws.Range("C24").Value = "Tweede der Landlieden"
$ws.Range("D24").Value = ""

$ws.Range("B25").Value = "#sylenus-mede-gedronken-hebbende"
$ws.Range("C25").Value = "Sylenus mede gedronken hebbende"
$ws.Range("D25").Value = ""

$ws.Range("B26").Value = "#maatzang"
$ws.Range("C26").Value = "Maatzang"
$ws.Range("D26").Value = ""

$ws.Range("B27").Value = "#een-uit-het-gevolg-van-ceres"
$ws.Range("C27").Value = "Een uit het gevolg van Ceres"
$ws.Range("D27").Value = ""

$ws.Range("B28").Value = "#mercurius"
$ws.Range("C28").Value = "Mercurius"
$ws.Range("D28").Value = ""

$ws.Range("B29").Value = "#maatzang-en-dichtkunde-te-zaamen"
$ws.Range("C29").Value = "Maatzang en Dichtkunde te zaamen"
$ws.Range("D29").Value = ""

$ws.Range("B30").Value = "#danskunst"
$ws.Range("C30").Value = "Danskunst"
$ws.Range("D30").Value = ""

$ws.Range("B31").Value = "#tweede-der-bacchanten-alleen"
$ws.Range("C31").Value = "Tweede der Bacchanten alleen"
$ws.Range("D31").Value = ""

$ws.Range("B32").Value = "#eerste-der-landtlieden"
$ws.Range("C32").Value = "Eerste der Landtlieden"
$ws.Range("D32").Value = ""

$ws.Range("B33").Value = "#gevolg-van-ceres"
$ws.Range("C33").Value = "Gevolg van Ceres"
$ws.Range("D33").Value = ""

$ws.Range("B34").Value = "#eene-der-bacchanten"
$ws.Range("C34").Value = "Eene der Bacchanten"
$ws.Range("D34").Value = ""

$ws.Range("B35").Value = "#eene-der-menaden"
$ws.Range("C35").Value = "Eene der Menaden"
$ws.Range("D35").Value = ""
